# Auto update: 2025-12-06 00:21:02
# Updates the daily semiconductor-stock decision table:
#  - refreshes the report date from 2025-12-05 to 2025-12-06
#  - refreshes price/RSI/score metrics for each ticker
#  - QCOM and AMD swap table rows (AMD now row 4, QCOM now row 6)
#  - ASML's judgement (column M) changes from the "stay out" to the
#    "buy-watch" message

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateCell {
    param($cell, [string]$value)
    # Force the cell to stay a plain text value instead of letting Excel
    # auto-convert the yyyy-mm-dd-looking string into a date serial.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- Row 2: ASML ----
Set-DateCell $ws.Range("A2") "2025-12-06"
$ws.Range("D2").Value = 1122.4
$ws.Range("E2").Value = 66.3
$ws.Range("F2").Value = 5.89
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 66
$ws.Range("K2").Value = 63.9
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 51.54219175917372

# ---- Row 3: TSM ----
Set-DateCell $ws.Range("A3") "2025-12-06"
$ws.Range("D3").Value = 299.95
$ws.Range("E3").Value = 63.4
$ws.Range("F3").Value = 2.89
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 58.7
$ws.Range("N3").Value = 51.54219175917372

# ---- Row 4: was QCOM, now AMD ----
Set-DateCell $ws.Range("A4") "2025-12-06"
$ws.Range("B4").Value = "Advanced Micro Devices, Inc."
$ws.Range("C4").Value = "AMD"
$ws.Range("D4").Value = 222.98
$ws.Range("E4").Value = 37.1
$ws.Range("F4").Value = 2.51
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 76
$ws.Range("I4").Value = 86
$ws.Range("J4").Value = 73
$ws.Range("K4").Value = 55.9
$ws.Range("N4").Value = 51.54219175917372

# ---- Row 5: NVDA ----
Set-DateCell $ws.Range("A5") "2025-12-06"
$ws.Range("D5").Value = 182.74
$ws.Range("E5").Value = 42
$ws.Range("F5").Value = 3.25
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 66
$ws.Range("J5").Value = 66
$ws.Range("K5").Value = 50.9
$ws.Range("N5").Value = 51.54219175917372

# ---- Row 6: was AMD, now QCOM ----
Set-DateCell $ws.Range("A6") "2025-12-06"
$ws.Range("B6").Value = "QUALCOMM Incorporated"
$ws.Range("C6").Value = "QCOM"
$ws.Range("D6").Value = 175.78
$ws.Range("E6").Value = 53.6
$ws.Range("F6").Value = 5.11
$ws.Range("G6").Value = 60
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 43
$ws.Range("K6").Value = 49.5
$ws.Range("N6").Value = 51.54219175917372

$wb.Save()
